# chore: update Sheets via scheduled runner
# Refreshes market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets with the
# latest scrape, row by row, leaving every other cell untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5635.6665
$ws.Range("I113").Value = 2287
$ws.Range("J113").Value = 7310
$ws.Range("K113").Value = 2287
$ws.Range("L113").Value = 7310
$ws.Range("M113").Value = 967
$ws.Range("N113").Value = -13818

$ws.Range("H116").Value = 4654.0557
$ws.Range("I116").Value = 3708
$ws.Range("K116").Value = 3708
$ws.Range("M116").Value = -266

$ws.Range("H127").Value = 2057.1428
$ws.Range("I127").Value = 1666.6666
$ws.Range("J127").Value = 4400
$ws.Range("K127").Value = 4999.9998
$ws.Range("L127").Value = 13200
$ws.Range("M127").Value = -39.9997999999996
$ws.Range("N127").Value = -23120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 38000
$ws.Range("J9").Value = 38000
$ws.Range("L9").Value = 38000
$ws.Range("N9").Value = -38340

$ws.Range("H20").Value = 38000
$ws.Range("J20").Value = 38000
$ws.Range("L20").Value = 38000
$ws.Range("N20").Value = -38540

$ws.Range("H61").Value = 2993.8
$ws.Range("J61").Value = 2991.3333
$ws.Range("L61").Value = 2991.3333
$ws.Range("N61").Value = -3415.3333

$ws.Range("H88").Value = 1181.8667
$ws.Range("I88").Value = 457.16666
$ws.Range("K88").Value = 457.16666
$ws.Range("M88").Value = -51.16665999999998

$ws.Range("H91").Value = 1181.8667
$ws.Range("I91").Value = 457.16666
$ws.Range("K91").Value = 457.16666
$ws.Range("M91").Value = 946.83334

$ws.Range("H122").Value = 835644.7
$ws.Range("I122").Value = 2426.45
$ws.Range("J122").Value = 4168517.5
$ws.Range("K122").Value = 7279.349999999999
$ws.Range("L122").Value = 12505552.5
$ws.Range("M122").Value = -4829.349999999999
$ws.Range("N122").Value = -12510452.5

$ws.Range("H132").Value = 1811.381
$ws.Range("I132").Value = 1521
$ws.Range("J132").Value = 2629.7273
$ws.Range("K132").Value = 4563
$ws.Range("L132").Value = 7889.1819
$ws.Range("M132").Value = -2033
$ws.Range("N132").Value = -12949.1819

$ws.Range("H136").Value = 2993.8
$ws.Range("J136").Value = 2991.3333
$ws.Range("L136").Value = 8973.999899999999
$ws.Range("N136").Value = -14073.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4355020
$ws.Range("I86").Value = 5007808
$ws.Range("K86").Value = 5007808
$ws.Range("M86").Value = -5006685

$ws.Range("H89").Value = 4355020
$ws.Range("I89").Value = 5007808
$ws.Range("K89").Value = 25039040
$ws.Range("M89").Value = -25033424

$ws.Range("H107").Value = 7144132.5
$ws.Range("I107").Value = 7937703
$ws.Range("K107").Value = 7937703
$ws.Range("M107").Value = -7935783

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37186.848
$ws.Range("I31").Value = 1466.3077
$ws.Range("K31").Value = 1466.3077
$ws.Range("M31").Value = -1171.3077

$ws.Range("H34").Value = 37186.848
$ws.Range("I34").Value = 1466.3077
$ws.Range("K34").Value = 1466.3077
$ws.Range("M34").Value = -1264.3077

$ws.Range("H93").Value = 47733
$ws.Range("J93").Value = 99799
$ws.Range("L93").Value = 99799
$ws.Range("N93").Value = -103543

$ws.Range("H99").Value = 3508.15
$ws.Range("I99").Value = 3214.2856
$ws.Range("J99").Value = 3666.3845
$ws.Range("K99").Value = 3214.2856
$ws.Range("L99").Value = 3666.3845
$ws.Range("M99").Value = -1716.2856
$ws.Range("N99").Value = -6662.3845

$ws.Range("H105").Value = 3109.111
$ws.Range("I105").Value = 3109.111
$ws.Range("K105").Value = 3109.111
$ws.Range("M105").Value = -1362.111

$ws.Range("H126").Value = 3508.15
$ws.Range("I126").Value = 3214.2856
$ws.Range("J126").Value = 3666.3845
$ws.Range("K126").Value = 9642.856800000001
$ws.Range("L126").Value = 10999.1535
$ws.Range("M126").Value = -7172.856800000001
$ws.Range("N126").Value = -15939.1535

$ws.Range("H134").Value = 3436.1724
$ws.Range("I134").Value = 2753.2778
$ws.Range("J134").Value = 4553.636
$ws.Range("K134").Value = 8259.8334
$ws.Range("L134").Value = 13660.908
$ws.Range("M134").Value = -5724.8334
$ws.Range("N134").Value = -18730.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1665
$ws.Range("I3").Value = 1558
$ws.Range("K3").Value = 4674
$ws.Range("M3").Value = -4562

$ws.Range("H80").Value = 2790.6
$ws.Range("I80").Value = 2674.8
$ws.Range("J80").Value = 2906.4
$ws.Range("K80").Value = 8024.400000000001
$ws.Range("L80").Value = 8719.200000000001
$ws.Range("M80").Value = -7088.400000000001
$ws.Range("N80").Value = -10591.2

$ws.Range("H83").Value = 2790.6
$ws.Range("I83").Value = 2674.8
$ws.Range("J83").Value = 2906.4
$ws.Range("K83").Value = 24073.2
$ws.Range("L83").Value = 26157.6
$ws.Range("M83").Value = -19393.2
$ws.Range("N83").Value = -35517.60000000001

$ws.Range("H121").Value = 799.5454999999999
$ws.Range("J121").Value = 1353
$ws.Range("L121").Value = 4059
$ws.Range("N121").Value = -6679

$ws.Range("H131").Value = 15435149
$ws.Range("I131").Value = 6411840.5
$ws.Range("J131").Value = 23813936
$ws.Range("K131").Value = 19235521.5
$ws.Range("L131").Value = 71441808
$ws.Range("M131").Value = -19230481.5
$ws.Range("N131").Value = -71451888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14292947
$ws.Range("I70").Value = 16674273
$ws.Range("J70").Value = 4993.5
$ws.Range("K70").Value = 16674273
$ws.Range("L70").Value = 4993.5
$ws.Range("M70").Value = -16674003
$ws.Range("N70").Value = -5533.5

$ws.Range("H73").Value = 14292947
$ws.Range("I73").Value = 16674273
$ws.Range("J73").Value = 4993.5
$ws.Range("K73").Value = 16674273
$ws.Range("L73").Value = 4993.5
$ws.Range("M73").Value = -16673337
$ws.Range("N73").Value = -6865.5

$ws.Range("H122").Value = 400597.4
$ws.Range("I122").Value = 539571.5
$ws.Range("K122").Value = 1618714.5
$ws.Range("M122").Value = -1616264.5

$ws.Range("H132").Value = 3063.7188
$ws.Range("I132").Value = 2545.625
$ws.Range("J132").Value = 4618
$ws.Range("K132").Value = 7636.875
$ws.Range("L132").Value = 13854
$ws.Range("M132").Value = -5106.875
$ws.Range("N132").Value = -18914

$ws.Range("H135").Value = 57000
$ws.Range("J135").Value = 57000
$ws.Range("L135").Value = 57000
$ws.Range("N135").Value = -67140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4436.154
$ws.Range("I7").Value = 1711.1666
$ws.Range("J7").Value = 6771.857
$ws.Range("K7").Value = 1711.1666
$ws.Range("L7").Value = 6771.857
$ws.Range("M7").Value = -1599.1666
$ws.Range("N7").Value = -6995.857

$ws.Range("H40").Value = 7234.353
$ws.Range("I40").Value = 5771.636
$ws.Range("J40").Value = 9916
$ws.Range("K40").Value = 5771.636
$ws.Range("L40").Value = 9916
$ws.Range("M40").Value = -5635.636
$ws.Range("N40").Value = -10188

$ws.Range("H46").Value = 4768.875
$ws.Range("I46").Value = 4033.3333
$ws.Range("K46").Value = 4033.3333
$ws.Range("M46").Value = -3845.3333

$ws.Range("H55").Value = 1781.0605
$ws.Range("I55").Value = 1578.8948
$ws.Range("K55").Value = 1578.8948
$ws.Range("M55").Value = -1405.8948

$ws.Range("H92").Value = 65000
$ws.Range("J92").Value = 65000
$ws.Range("L92").Value = 65000
$ws.Range("N92").Value = -69992

$ws.Range("H122").Value = 4261.25
$ws.Range("I122").Value = 2625.158
$ws.Range("K122").Value = 7875.474
$ws.Range("M122").Value = -5425.474

$ws.Range("H126").Value = 4436.154
$ws.Range("I126").Value = 1711.1666
$ws.Range("J126").Value = 6771.857
$ws.Range("K126").Value = 5133.4998
$ws.Range("L126").Value = 20315.571
$ws.Range("M126").Value = -2663.4998
$ws.Range("N126").Value = -25255.571

$ws.Range("H136").Value = 92065.30499999999
$ws.Range("I136").Value = 137460.2
$ws.Range("J136").Value = 6949.875
$ws.Range("K136").Value = 412380.6
$ws.Range("L136").Value = 20849.625
$ws.Range("M136").Value = -409830.6
$ws.Range("N136").Value = -25949.625
